$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.339, 0.094, 0.123, 107, 14, 83),
    @(0.135, -0.043, 0.322, -162, 13, -56),
    @(-0.092, -0.08500000000000001, 0.249, 21, -32, 65),
    @(0.256, -0.101, 0.435, 157, -80, -67),
    @(0.161, 0.095, -0.241, 92, 1, -55),
    @(0.208, 0.083, 0.08400000000000001, 90, 0, 157),
    @(0.104, -0.042, 0.333, -87, -82, 88),
    @(-0.003, -0.351, -0.002, 114, -25, 25),
    @(0.333, -0.049, 0.184, 153, -29, 113),
    @(0.04, 0.126, 0.408, 76, 21, 168)
)

for ($r = 1; $r -le 10; $r++) {
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $data[$r - 1][$c - 1]
    }
}
